$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3.83
    3  = 0.24
    4  = 0.91
    5  = 0.16
    6  = 5.12
    7  = 0.26
    8  = 0.25
    9  = 0.19
    10 = 0.11
    12 = 4.27
    13 = 0.86
    14 = 14.53
    15 = 1.64
    16 = 0.3
    17 = 4.23
    18 = 0.17
    19 = 0.12
    20 = 0.15
    22 = 0.26
    23 = 0.77
    24 = 0.8
    25 = 0.29
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 4).Value = $updates[$row]
}
